# Update the "想去人数" (interested-attendee count) figures that changed
# between site regenerations.
#
# Sheet "展览"  (sheet1 / index 1): rows 4-7 in column F
# Sheet "全部类型" (sheet4 / index 4): rows 4, 6, 7, 8 in column F
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 1771
$wsExhibit.Range("F5").Value = 788
$wsExhibit.Range("F6").Value = 257
$wsExhibit.Range("F7").Value = 206

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1771
$wsAll.Range("F6").Value = 788
$wsAll.Range("F7").Value = 257
$wsAll.Range("F8").Value = 206
